$wb = $excel.ActiveWorkbook

# --- Move the live selection on currency_conversions to A3 (this sheet is
#     no longer the active tab afterwards, but keeps this resting selection).
$wsCur = $wb.Worksheets.Item("currency_conversions")
$wsCur.Activate()
$wsCur.Range("A3").Select()

# --- Add the new "currency_movements" sheet as the last tab, right after
#     currency_conversions.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "currency_movements"

# --- Header row (bold) ---
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "buy_date"
$ws.Range("C1").Value = "amount"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "comment"
$ws.Range("A1:E1").Font.Bold = $true

# --- Data rows ---
$fmt = "yyyy\-mm\-dd;@"

$ws.Range("A2").NumberFormat = $fmt
$ws.Range("A2").Value = Get-Date -Year 2022 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("B2").NumberFormat = $fmt
$ws.Range("B2").Value = Get-Date -Year 2022 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "USD"

$ws.Range("A3").NumberFormat = $fmt
$ws.Range("A3").Value = Get-Date -Year 2022 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("B3").NumberFormat = $fmt
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = "USD"

$ws.Range("A4").NumberFormat = $fmt
$ws.Range("A4").Value = Get-Date -Year 2022 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "EUR"

$ws.Range("A5").NumberFormat = $fmt
$ws.Range("A5").Value = Get-Date -Year 2022 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("C5").Value = -1
$ws.Range("D5").Value = "EUR"

# --- New sheet becomes the active tab/selection, as it would after a user
#     adds a sheet and starts filling it in.
$ws.Activate()
$ws.Range("A1").Select()
